# Regenerate orders with updated distance/sizes.
# The workbook encodes experiment trial metadata as strings such as
# "Face13_D64_S30" / "Face13_D64_S30_l.png" etc. This edit renames the
# Distance codes and one of the Size codes throughout the whole sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (S25 and S20 are left untouched.)
#
# Doing this as whole-sheet text replacements reproduces the same effect
# as regenerating the order file with the new codes, without disturbing
# any other cell, numeric value, or formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Order matters: none of the new codes collide with any of the old codes,
# but replace the longer/more specific distance tokens before the size
# token just to keep things unambiguous and easy to follow.
$used.Replace("D64", "D69")
$used.Replace("D80", "D86")
$used.Replace("D51", "D55")
$used.Replace("S30", "S31")
